# "(completely) new version (stable)"
# Update a handful of Russian UI-label strings on the translation sheet
# and move the active selection to D48 (matching the author's saved view).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# auth.text (ru + en columns both held the Russian placeholder copy)
$ws.Range("C37").Value = "Первым делом необходимо авторизоваться. <br>Получить пароль можно у"
$ws.Range("D37").Value = "Первым делом необходимо авторизоваться. <br>Получить пароль можно у"

# geo2 ("Прочие примечания" -> "Примечания к расположению")
$ws.Range("C46").Value = "Примечания к расположению"

# ev.title ("Данные пробы" -> "Сбор материала")
$ws.Range("C49").Value = "Сбор материала"

# ev.rem ("Комментарии" -> "Примечания к сбору материала")
$ws.Range("C54").Value = "Примечания к сбору материала"

# Match the saved selection from the commit (cell D48, row 48's Coordinates: value).
$ws.Range("D48").Select()
